$d = $word.ActiveDocument

$replacements = @(
    @("2025-02-03 Monday", "2025-02-04 Tuesday"),
    @("72÷3=", "16÷7="),
    @("26÷4=", "52÷3="),
    @("11÷3=", "93÷4="),
    @("88÷3=", "29÷7="),
    @("78÷8=", "47÷3="),
    @("56÷4=", "50÷5="),
    @("77÷2=", "26÷6="),
    @("95÷8=", "52÷7="),
    @("20÷5=", "14÷6="),
    @("58÷9=", "68÷7="),
    @("10÷4=", "38÷7="),
    @("76÷3=", "48÷4="),
    @("32÷4=", "10÷2="),
    @("78÷7=", "13÷2="),
    @("30÷2=", "34÷9="),
    @("48÷3=", "85÷5="),
    @("70÷8=", "28÷2="),
    @("82÷7=", "15÷9="),
    @("94÷2=", "77÷5="),
    @("26÷2=", "92÷7="),
    @("88÷8=", "72÷9="),
    @("19÷5=", "86÷3="),
    @("85÷3=", "38÷6="),
    @("76÷8=", "61÷4="),
    @("64÷5=", "82÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
